# "Generate Report for Archive"
#
# The localization status report is regenerated:
#  - 4fa03f76-...md and f3f2c32c-...md move from "Ready for handoff" to
#    "In Translation".
#  - f3f2c32c-...md and 6be85808-...md swap positions in the per-language
#    tables (and the Overview sheet), because f3f2c32c is now ahead of
#    6be85808 in the list.
#  - The "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" and
#    "Latest Handoff File" values follow the files that moved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 (4fa03f76): status changes to "In Translation"
$ws.Cells.Item(3, 5).Value = "In Translation"
$ws.Cells.Item(3, 6).Value = "In Translation"

# Row 4 used to be 6be85808, now becomes f3f2c32c ("In Translation")
$ws.Cells.Item(4, 1).Value = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
$ws.Cells.Item(4, 2).Value = "e2e\f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
$ws.Cells.Item(4, 5).Value = "In Translation"
$ws.Cells.Item(4, 6).Value = "In Translation"
$ws.Cells.Item(4, 7).Value = "2016-08-20 14:44:31"

# Row 5 used to be f3f2c32c, now becomes 6be85808 ("Ready for handoff")
$ws.Cells.Item(5, 1).Value = "6be85808-a57d-4e45-97a4-8d73994503e5.md"
$ws.Cells.Item(5, 2).Value = "e2e\6be85808-a57d-4e45-97a4-8d73994503e5.md"
$ws.Cells.Item(5, 5).Value = "Ready for handoff"
$ws.Cells.Item(5, 6).Value = "Ready for handoff"
$ws.Cells.Item(5, 7).Value = "2016-08-20 14:43:31"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
    }
    elseif ($addr -eq '$B$5') {
        $hl.TextToDisplay = "e2e\6be85808-a57d-4e45-97a4-8d73994503e5.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3 (4fa03f76): status changes to "In Translation"
$ws.Cells.Item(3, 3).Value = "In Translation"

# Row 4 used to be 6be85808, now becomes f3f2c32c ("In Translation")
$ws.Cells.Item(4, 1).Value = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
$ws.Cells.Item(4, 3).Value = "In Translation"
$ws.Cells.Item(4, 7).Value = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.20d26e755cd9a3cf4ae816105b88ea1895c3e196.zh-cn.xlf"
$ws.Cells.Item(4, 8).Value = "2016-08-20 14:44:27"

# Row 5 used to be f3f2c32c, now becomes 6be85808 ("Ready for handoff")
$ws.Cells.Item(5, 1).Value = "6be85808-a57d-4e45-97a4-8d73994503e5.md"
$ws.Cells.Item(5, 7).Value = "6be85808-a57d-4e45-97a4-8d73994503e5.15b5a30a64c2df2a913b6b8b3db4b2b15f0bed3c.zh-cn.xlf"
$ws.Cells.Item(5, 8).Value = "2016-08-20 14:43:27"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') {
        $hl.TextToDisplay = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
    }
    elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "6be85808-a57d-4e45-97a4-8d73994503e5.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 3 (4fa03f76): status changes to "In Translation"
$ws.Cells.Item(3, 3).Value = "In Translation"

# Row 4 used to be 6be85808, now becomes f3f2c32c ("In Translation")
$ws.Cells.Item(4, 1).Value = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
$ws.Cells.Item(4, 3).Value = "In Translation"
$ws.Cells.Item(4, 7).Value = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.20d26e755cd9a3cf4ae816105b88ea1895c3e196.de-de.xlf"
$ws.Cells.Item(4, 8).Value = "2016-08-20 14:44:31"

# Row 5 used to be f3f2c32c, now becomes 6be85808 ("Ready for handoff")
$ws.Cells.Item(5, 1).Value = "6be85808-a57d-4e45-97a4-8d73994503e5.md"
$ws.Cells.Item(5, 7).Value = "6be85808-a57d-4e45-97a4-8d73994503e5.15b5a30a64c2df2a913b6b8b3db4b2b15f0bed3c.de-de.xlf"
$ws.Cells.Item(5, 8).Value = "2016-08-20 14:43:31"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') {
        $hl.TextToDisplay = "f3f2c32c-c6c0-43e2-8590-13f07f41e0af.md"
    }
    elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "6be85808-a57d-4e45-97a4-8d73994503e5.md"
    }
}
